$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Copy the formatting of the existing header cell (H1) onto the two new
# header cells so they pick up the same bold/border/centered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2 is a special case (does not follow the I=1 / J=H pattern)
$ws.Cells.Item(2, 9).Value = 4
$ws.Cells.Item(2, 10).Value = 6

# Rows 3 through 40: I = 1, J = existing H value
for ($r = 3; $r -le 40; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
